$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1918367346938775
$ws.Range("C2").Value = 0.5571428571428572
$ws.Range("J2").Value = 0.01224489795918367
$ws.Range("P2").Value = 0.1306122448979592
$ws.Range("S2").Value = 0.1081632653061225
$ws.Range("B3").Value = 0.01067615658362989
$ws.Range("C3").Value = 0.02135231316725979
$ws.Range("J3").Value = 0.01779359430604982
$ws.Range("P3").Value = 0.7402135231316725
$ws.Range("S3").Value = 0.2099644128113879
$ws.Range("J4").Value = 0.0303030303030303
$ws.Range("P4").Value = 0.6868686868686869
$ws.Range("S4").Value = 0.2828282828282828
$ws.Range("B6").Value = 0.05665024630541872
$ws.Range("D6").Value = 0.01724137931034483
$ws.Range("F6").Value = 0.06896551724137931
$ws.Range("J6").Value = 0.2438423645320197
$ws.Range("O6").Value = 0.01231527093596059
$ws.Range("Q6").Value = 0.1748768472906404
$ws.Range("R6").Value = 0.07389162561576355
$ws.Range("S6").Value = 0.3522167487684729
$ws.Range("B7").Value = 0.09785202863961814
$ws.Range("D7").Value = 0.02625298329355609
$ws.Range("E7").Value = 0.00477326968973747
$ws.Range("F7").Value = 0.05011933174224344
$ws.Range("J7").Value = 0.1264916467780429
$ws.Range("O7").Value = 0.01193317422434368
$ws.Range("Q7").Value = 0.2219570405727924
$ws.Range("R7").Value = 0.03818615751789976
$ws.Range("S7").Value = 0.4224343675417661
$ws.Range("B8").Value = 0.07367280606717226
$ws.Range("D8").Value = 0.02491874322860238
$ws.Range("F8").Value = 0.06392199349945829
$ws.Range("J8").Value = 0.1115926327193933
$ws.Range("O8").Value = 0.01841820151679307
$ws.Range("Q8").Value = 0.2047670639219935
$ws.Range("R8").Value = 0.07583965330444203
$ws.Range("S8").Value = 0.4268689057421452
$ws.Range("B9").Value = 0.08808290155440414
$ws.Range("D9").Value = 0.02072538860103627
$ws.Range("F9").Value = 0.05440414507772021
$ws.Range("J9").Value = 0.1088082901554404
$ws.Range("O9").Value = 0.02072538860103627
$ws.Range("Q9").Value = 0.1735751295336788
$ws.Range("R9").Value = 0.08808290155440414
$ws.Range("S9").Value = 0.4455958549222798
$ws.Range("B10").Value = 0.09297433740008414
$ws.Range("D10").Value = 0.02271771140092554
$ws.Range("E10").Value = 0.0008413967185527977
$ws.Range("F10").Value = 0.06604964240639462
$ws.Range("J10").Value = 0.1127471602860749
$ws.Range("O10").Value = 0.01304164913756836
$ws.Range("Q10").Value = 0.2254943205721498
$ws.Range("R10").Value = 0.07867059318468658
$ws.Range("S10").Value = 0.3874631888935633
$ws.Range("G11").Value = 0.1257575757575758
$ws.Range("J11").Value = 0.0893939393939394
$ws.Range("K11").Value = 0.1893939393939394
$ws.Range("L11").Value = 0.5863636363636363
$ws.Range("S11").Value = 0.00909090909090909
$ws.Range("G12").Value = 0.7386934673366834
$ws.Range("J12").Value = 0.185929648241206
$ws.Range("K12").Value = 0.005025125628140704
$ws.Range("L12").Value = 0.02261306532663317
$ws.Range("S12").Value = 0.04773869346733668
$ws.Range("G13").Value = 0.5888888888888889
$ws.Range("J13").Value = 0.3444444444444444
$ws.Range("S13").Value = 0.06666666666666667
$ws.Range("F15").Value = 0.01541850220264317
$ws.Range("H15").Value = 0.1696035242290749
$ws.Range("I15").Value = 0.08590308370044053
$ws.Range("J15").Value = 0.3215859030837004
$ws.Range("K15").Value = 0.08590308370044053
$ws.Range("M15").Value = 0.00881057268722467
$ws.Range("O15").Value = 0.09251101321585903
$ws.Range("S15").Value = 0.2202643171806167
$ws.Range("F16").Value = 0.0182370820668693
$ws.Range("H16").Value = 0.1762917933130699
$ws.Range("I16").Value = 0.0911854103343465
$ws.Range("J16").Value = 0.39209726443769
$ws.Range("K16").Value = 0.1337386018237082
$ws.Range("M16").Value = 0.0121580547112462
$ws.Range("N16").Value = 0.00303951367781155
$ws.Range("O16").Value = 0.060790273556231
$ws.Range("S16").Value = 0.1124620060790274
$ws.Range("F17").Value = 0.01888772298006296
$ws.Range("H17").Value = 0.1951731374606506
$ws.Range("I17").Value = 0.09548793284365163
$ws.Range("J17").Value = 0.3767051416579223
$ws.Range("K17").Value = 0.1017838405036726
$ws.Range("M17").Value = 0.0167890870933893
$ws.Range("N17").Value = 0.001049317943336831
$ws.Range("O17").Value = 0.0608604407135362
$ws.Range("S17").Value = 0.1332633788037775
$ws.Range("F18").Value = 0.01780415430267062
$ws.Range("H18").Value = 0.172106824925816
$ws.Range("I18").Value = 0.09495548961424333
$ws.Range("J18").Value = 0.4243323442136498
$ws.Range("K18").Value = 0.09495548961424333
$ws.Range("M18").Value = 0.008902077151335312
$ws.Range("O18").Value = 0.05044510385756677
$ws.Range("S18").Value = 0.1364985163204748
$ws.Range("F19").Value = 0.01614763552479815
$ws.Range("H19").Value = 0.2106881968473664
$ws.Range("I19").Value = 0.07650903498654364
$ws.Range("J19").Value = 0.340638216070742
$ws.Range("K19").Value = 0.1222606689734717
$ws.Range("M19").Value = 0.02691272587466359
$ws.Range("N19").Value = 0.001537870049980777
$ws.Range("O19").Value = 0.0722798923490965
$ws.Range("S19").Value = 0.1330257593233372
